# edit.ps1 - apply the FLARIUM.docx changes described by the diff:
#  1. Remove the stray "_GoBack" bookmark around "Discussion Owner".
#  2. Turn "Close the discussion" into "Update(title, closed_at) the
#     discussion" (highlighted, with closed_at spell-check wrapped) and
#     highlight the two sub-bullets that already follow it.
#  3. Insert a new "Check if the discussi[on is open.]" sub-bullet (with a
#     fresh "_GoBack" bookmark inside it) plus a whole new "Delete the
#     discussion" block (with its own two sub-bullets) before "Channels".
#  4. Remove the yellow highlight from the "Pagination" bullet.
#  5. Add a lastRenderedPageBreak marker to the final paragraph.

$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------
# 1. Drop the leftover _GoBack bookmark on the "Discussion Owner" line.
# ---------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------
# 2. Locate the "Close the discussion" bullet group by scanning text
#    (paragraph indices below this point do not move until step 3, since
#    every earlier edit is a 1-for-1 paragraph replacement).
# ---------------------------------------------------------------------
$closeIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq "Close the discussion") {
        $closeIdx = $i
        break
    }
}

$authIdx = $closeIdx + 1
$ownerIdx = $closeIdx + 2

# "Close the discussion" -> "Update(title, closed_at) the discussion"
$p = $d.Paragraphs($closeIdx)
$r = $p.Range.Duplicate
$xml = "<w:p $wNs>" +
       "<w:pPr><w:pStyle w:val='Paragraphedeliste'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='4'/></w:numPr><w:rPr><w:highlight w:val='yellow'/></w:rPr></w:pPr>" +
       "<w:r><w:rPr><w:highlight w:val='yellow'/></w:rPr><w:t xml:space='preserve'>Update(title, </w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/>" +
       "<w:r><w:rPr><w:highlight w:val='yellow'/></w:rPr><w:t>closed_at</w:t></w:r>" +
       "<w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:rPr><w:highlight w:val='yellow'/></w:rPr><w:t>)</w:t></w:r>" +
       "<w:r><w:rPr><w:highlight w:val='yellow'/></w:rPr><w:t xml:space='preserve'> the discussion</w:t></w:r>" +
       "</w:p>"
$r.InsertXML($xml)

# "Authentication required" (sub-bullet) -> same text, now highlighted.
$p = $d.Paragraphs($authIdx)
$r = $p.Range.Duplicate
$xml = "<w:p $wNs>" +
       "<w:pPr><w:pStyle w:val='Paragraphedeliste'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='4'/></w:numPr><w:rPr><w:highlight w:val='yellow'/></w:rPr></w:pPr>" +
       "<w:r><w:rPr><w:highlight w:val='yellow'/></w:rPr><w:t>Authentication required</w:t></w:r>" +
       "</w:p>"
$r.InsertXML($xml)

# "Check if the user is the owner" (sub-bullet, 3 runs) -> same text, highlighted.
$p = $d.Paragraphs($ownerIdx)
$r = $p.Range.Duplicate
$xml = "<w:p $wNs>" +
       "<w:pPr><w:pStyle w:val='Paragraphedeliste'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='4'/></w:numPr><w:rPr><w:highlight w:val='yellow'/></w:rPr></w:pPr>" +
       "<w:r><w:rPr><w:highlight w:val='yellow'/></w:rPr><w:t xml:space='preserve'>Check if the user is </w:t></w:r>" +
       "<w:r><w:rPr><w:highlight w:val='yellow'/></w:rPr><w:t xml:space='preserve'>the </w:t></w:r>" +
       "<w:r><w:rPr><w:highlight w:val='yellow'/></w:rPr><w:t>owner</w:t></w:r>" +
       "</w:p>"
$r.InsertXML($xml)

# ---------------------------------------------------------------------
# 3. De-highlight the "Pagination" bullet (pPr mark + run both lose the
#    yellow highlight run property). Do this before the insertion below
#    so paragraph indices are still untouched by new paragraphs.
# ---------------------------------------------------------------------
$pagIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq "Pagination") {
        $pagIdx = $i
        break
    }
}
$p = $d.Paragraphs($pagIdx)
$r = $p.Range.Duplicate
$xml = "<w:p $wNs>" +
       "<w:pPr><w:pStyle w:val='Paragraphedeliste'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='4'/></w:numPr></w:pPr>" +
       "<w:r><w:t>Pagination</w:t></w:r>" +
       "</w:p>"
$r.InsertXML($xml)

# ---------------------------------------------------------------------
# 4. lastRenderedPageBreak on the very last paragraph of the document.
# ---------------------------------------------------------------------
$lastIdx = $d.Paragraphs.Count
$p = $d.Paragraphs($lastIdx)
$text = $p.Range.Text.TrimEnd()
$r = $p.Range.Duplicate
$xml = "<w:p $wNs>" +
       "<w:pPr><w:pStyle w:val='Paragraphedeliste'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='6'/></w:numPr></w:pPr>" +
       "<w:r><w:lastRenderedPageBreak/><w:t>$text</w:t></w:r>" +
       "</w:p>"
$r.InsertXML($xml)

# ---------------------------------------------------------------------
# 5. Insert the six new paragraphs right before "Channels": a new
#    "Check if the discussi[on is open.]" sub-bullet (with its own fresh
#    _GoBack bookmark), a spacer, a whole new "Delete the discussion"
#    block (+2 sub-bullets) and a second spacer.
#
#    NOTE: InsertXML on a *collapsed* (zero-length) Range behaves like a
#    splice into the surrounding paragraph (it can clobber the paragraph
#    that owns that boundary position) rather than a clean block insert.
#    The reliable primitive here is Range.InsertParagraphBefore(), which
#    always creates a genuine new empty paragraph immediately before the
#    target without touching the target's own content; we reserve six
#    empty slots that way and then fill each one via a whole-paragraph
#    InsertXML replace (safe once the paragraph is its own empty node).
# ---------------------------------------------------------------------
$chanIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq "Channels") {
        $chanIdx = $i
        break
    }
}

$newXmls = @(
    ("<w:p $wNs>" +
     "<w:pPr><w:pStyle w:val='Paragraphedeliste'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='4'/></w:numPr></w:pPr>" +
     "<w:r><w:t>Check if the discussi</w:t></w:r>" +
     "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
     "<w:r><w:t>on is open.</w:t></w:r>" +
     "</w:p>"),
    ("<w:p $wNs>" +
     "<w:pPr><w:ind w:left='1080'/><w:rPr><w:highlight w:val='yellow'/></w:rPr></w:pPr>" +
     "</w:p>"),
    ("<w:p $wNs>" +
     "<w:pPr><w:pStyle w:val='Paragraphedeliste'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='4'/></w:numPr><w:rPr><w:highlight w:val='yellow'/></w:rPr></w:pPr>" +
     "<w:r><w:rPr><w:highlight w:val='yellow'/></w:rPr><w:t>Delete the discussion</w:t></w:r>" +
     "</w:p>"),
    ("<w:p $wNs>" +
     "<w:pPr><w:pStyle w:val='Paragraphedeliste'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='4'/></w:numPr><w:rPr><w:highlight w:val='yellow'/></w:rPr></w:pPr>" +
     "<w:r><w:rPr><w:highlight w:val='yellow'/></w:rPr><w:t>Authentication required</w:t></w:r>" +
     "</w:p>"),
    ("<w:p $wNs>" +
     "<w:pPr><w:pStyle w:val='Paragraphedeliste'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='4'/></w:numPr><w:rPr><w:highlight w:val='yellow'/></w:rPr></w:pPr>" +
     "<w:r><w:rPr><w:highlight w:val='yellow'/></w:rPr><w:t>Check if the user is the owner</w:t></w:r>" +
     "</w:p>"),
    ("<w:p $wNs>" +
     "<w:pPr><w:pStyle w:val='Paragraphedeliste'/><w:ind w:left='1440'/><w:rPr><w:highlight w:val='yellow'/></w:rPr></w:pPr>" +
     "</w:p>")
)

$target = $d.Paragraphs($chanIdx)
for ($k = 0; $k -lt $newXmls.Count; $k++) {
    $target.Range.InsertParagraphBefore()
}
for ($k = 0; $k -lt $newXmls.Count; $k++) {
    $slot = $chanIdx + $k
    $r = $d.Paragraphs($slot).Range.Duplicate
    $r.InsertXML($newXmls[$k])
}
